$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 165 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "2021-10-13") {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = "2021-10-28"
        $cell.Style = $origStyle
    }
}
